$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Wipe the previous layout (content + formatting) for the area we rebuild.
# ---------------------------------------------------------------------------
$ws.Range("B7:E35").Clear()

# ---------------------------------------------------------------------------
# 2. Title block (unchanged content, same styles as before).
# ---------------------------------------------------------------------------
$c = $ws.Range("D7")
$c.Value2 = "EMPRESA DE OBRAS SANITARIAS DE PASTO EMPOPASTO "
$c.Font.Bold = $true
$c.Font.Color = 9527094
$c.HorizontalAlignment = -4108
$c.VerticalAlignment = -4108

$c = $ws.Range("B9")
$c.Value2 = "Día y hora de Impresión"
$c.Font.Bold = $true

$c = $ws.Range("D11")
$c.Value2 = "REGISTRO INTERNO DE PROVEEDORES"
$c.Font.Bold = $true
$c.Font.Size = 18
$c.Font.Color = 9527094
$c.HorizontalAlignment = -4108
$c.VerticalAlignment = -4108
$ws.Rows.Item(11).RowHeight = 23.25

# ---------------------------------------------------------------------------
# 3. Field labels (column B) with a blank, left-aligned text cell in column E
#    next to each one (the new style used there: text number format + left
#    alignment).
# ---------------------------------------------------------------------------
$labels = @(
    "CC/NIT",
    "Tipo de identificación",
    "Tipo de entidad",
    "Nombre o razón social",
    "Ciudad",
    "Dirección",
    "Email",
    "Teléfono",
    "Fax",
    "Celular",
    "Nombre contacto",
    "Apellido contacto",
    "Tarjeta profesional",
    "Rut",
    "Ponderado"
)

$row = 13
foreach ($lab in $labels) {
    $bCell = $ws.Cells.Item($row, 2)
    $bCell.Value2 = $lab
    $bCell.Font.Bold = $true

    $eCell = $ws.Cells.Item($row, 5)
    $eCell.NumberFormat = "@"
    $eCell.HorizontalAlignment = -4131

    $row = $row + 1
}

# ---------------------------------------------------------------------------
# 4. "ACTIVIDADES REGISTRADAS" sub-title + small table header row.
# ---------------------------------------------------------------------------
$c = $ws.Range("D29")
$c.Value2 = "ACTIVIDADES REGISTRADAS"
$c.Font.Bold = $true
$c.Font.Size = 18
$c.Font.Color = 9527094
$c.HorizontalAlignment = -4108
$c.VerticalAlignment = -4108
$ws.Rows.Item(29).RowHeight = 23.25

$c = $ws.Range("D30")
$c.HorizontalAlignment = -4108
$c.VerticalAlignment = -4108

$ws.Range("B31").Value2 = "Sección"
$ws.Range("C31").Value2 = "División"

$c = $ws.Range("D31")
$c.Value2 = "Actividad"
$c.HorizontalAlignment = -4108
$c.VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# 5. Reposition the logo picture (same size, shifted right and down).
# ---------------------------------------------------------------------------
$shp = $ws.Shapes.Item(1)
$shp.Left = 2084295 / 12700.0
$shp.Top = 179293 / 12700.0
$shp.Width = 2867025 / 12700.0
$shp.Height = 819150 / 12700.0

# ---------------------------------------------------------------------------
# 6. View state: scroll so row 13 is at the top, selection on G6.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G6").Select()
